# Insert a new row above row 39. This pushes the existing rows 39-211
# down to 40-212, preserving their data (including formatting/styles).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Insert()

# The row that used to be at 39 is now at row 40; copy its contents into
# the newly-created (empty) row 39, then overwrite the date (column D)
# with the new value for this entry.
$ws.Range("A40:R40").Copy()
$ws.Range("A39:R39").PasteSpecial()

$ws.Range("D39").Value = 44575
